$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 160.33333
$ws.Range("I11").Value = 160.33333
$ws.Range("K11").Value = 160.33333
$ws.Range("M11").Value = -20.33332999999999
$ws.Range("H19").Value = 1328.875
$ws.Range("I19").Value = 1637.8
$ws.Range("J19").Value = 814
$ws.Range("K19").Value = 1637.8
$ws.Range("L19").Value = 814
$ws.Range("M19").Value = -1462.8
$ws.Range("N19").Value = -1164
$ws.Range("H33").Value = 6025.1113
$ws.Range("I33").Value = 6930.3335
$ws.Range("K33").Value = 6930.3335
$ws.Range("M33").Value = -6701.3335
$ws.Range("H62").Value = 13141.412
$ws.Range("J62").Value = 9658.846
$ws.Range("L62").Value = 9658.846
$ws.Range("N62").Value = -10906.846
$ws.Range("H65").Value = 13141.412
$ws.Range("J65").Value = 9658.846
$ws.Range("L65").Value = 48294.23
$ws.Range("N65").Value = -54534.23
$ws.Range("H70").Value = 2504.8572
$ws.Range("I70").Value = 2136
$ws.Range("K70").Value = 6408
$ws.Range("M70").Value = -6138
$ws.Range("H73").Value = 2504.8572
$ws.Range("I73").Value = 2136
$ws.Range("K73").Value = 6408
$ws.Range("M73").Value = -5472
$ws.Range("H74").Value = 5165.6665
$ws.Range("I74").Value = 5165.6665
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 5165.6665
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -4229.6665
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 5165.6665
$ws.Range("I77").Value = 5165.6665
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 25828.3325
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -21148.3325
$ws.Range("N77").ClearContents()
$ws.Range("H97").Value = 2086.0908
$ws.Range("J97").Value = 2612.125
$ws.Range("L97").Value = 7836.375
$ws.Range("N97").Value = -8828.375
$ws.Range("H132").Value = 8375.932000000001
$ws.Range("I132").Value = 7114.231
$ws.Range("K132").Value = 21342.693
$ws.Range("M132").Value = -18812.693
$ws.Range("H138").Value = 3831.0408
$ws.Range("I138").Value = 2202.1035
$ws.Range("J138").Value = 6193
$ws.Range("K138").Value = 6606.310500000001
$ws.Range("L138").Value = 18579
$ws.Range("M138").Value = -1466.310500000001
$ws.Range("N138").Value = -28859
$ws.Range("H141").Value = 4828.2
$ws.Range("I141").Value = 4828.2
$ws.Range("K141").Value = 14484.6
$ws.Range("M141").Value = -9304.599999999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5818.75
$ws.Range("I63").Value = 2183.3333
$ws.Range("K63").Value = 2183.3333
$ws.Range("M63").Value = -1497.3333
$ws.Range("H66").Value = 5818.75
$ws.Range("I66").Value = 2183.3333
$ws.Range("K66").Value = 10916.6665
$ws.Range("M66").Value = -7484.666499999999
$ws.Range("H88").Value = 2107.5715
$ws.Range("I88").Value = 1117.6666
$ws.Range("J88").Value = 2377.5454
$ws.Range("K88").Value = 1117.6666
$ws.Range("L88").Value = 2377.5454
$ws.Range("M88").Value = -711.6666
$ws.Range("N88").Value = -3189.5454
$ws.Range("H91").Value = 2107.5715
$ws.Range("I91").Value = 1117.6666
$ws.Range("J91").Value = 2377.5454
$ws.Range("K91").Value = 1117.6666
$ws.Range("L91").Value = 2377.5454
$ws.Range("M91").Value = 286.3334
$ws.Range("N91").Value = -5185.5454
$ws.Range("H122").Value = 1716.4736
$ws.Range("I122").Value = 1494.9412
$ws.Range("K122").Value = 4484.8236
$ws.Range("M122").Value = -2034.8236
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3863.3572
$ws.Range("I20").Value = 3698.8572
$ws.Range("J20").Value = 4027.8572
$ws.Range("K20").Value = 3698.8572
$ws.Range("L20").Value = 4027.8572
$ws.Range("M20").Value = -3451.8572
$ws.Range("N20").Value = -4521.8572
$ws.Range("H134").Value = 2485.2942
$ws.Range("J134").Value = 5197.25
$ws.Range("L134").Value = 15591.75
$ws.Range("N134").Value = -20661.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 981.8333
$ws.Range("I16").Value = 981.8333
$ws.Range("K16").Value = 981.8333
$ws.Range("M16").Value = -694.8333
$ws.Range("H31").Value = 2381.7144
$ws.Range("I31").Value = 2139.7144
$ws.Range("J31").Value = 3349.7144
$ws.Range("K31").Value = 2139.7144
$ws.Range("L31").Value = 3349.7144
$ws.Range("M31").Value = -1844.7144
$ws.Range("N31").Value = -3939.7144
$ws.Range("H34").Value = 2381.7144
$ws.Range("I34").Value = 2139.7144
$ws.Range("J34").Value = 3349.7144
$ws.Range("K34").Value = 2139.7144
$ws.Range("L34").Value = 3349.7144
$ws.Range("M34").Value = -1937.7144
$ws.Range("N34").Value = -3753.7144
$ws.Range("H113").Value = 981.8333
$ws.Range("I113").Value = 981.8333
$ws.Range("K113").Value = 981.8333
$ws.Range("M113").Value = 1188.1667
$ws.Range("H122").Value = 2352.375
$ws.Range("I122").Value = 2389.5334
$ws.Range("K122").Value = 7168.600199999999
$ws.Range("M122").Value = -4718.600199999999
$ws.Range("H134").Value = 75205.64
$ws.Range("I134").Value = 113008.89
$ws.Range("J134").Value = 7159.8
$ws.Range("K134").Value = 339026.67
$ws.Range("L134").Value = 21479.4
$ws.Range("M134").Value = -336491.67
$ws.Range("N134").Value = -26549.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1249.7693
$ws.Range("J34").Value = 3000
$ws.Range("L34").Value = 9000
$ws.Range("N34").Value = -9168
$ws.Range("H39").Value = 9660
$ws.Range("J39").Value = 9660
$ws.Range("L39").Value = 28980
$ws.Range("N39").Value = -29568
$ws.Range("H55").Value = 4062.6667
$ws.Range("I55").Value = 4062.6667
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 12188.0001
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -12011.0001
$ws.Range("N55").ClearContents()
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H140").Value = 2480.4443
$ws.Range("I140").Value = 2024.9375
$ws.Range("K140").Value = 6074.8125
$ws.Range("M140").Value = -894.8125
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14458.333
$ws.Range("I70").Value = 14687.75
$ws.Range("K70").Value = 14687.75
$ws.Range("M70").Value = -14417.75
$ws.Range("H73").Value = 14458.333
$ws.Range("I73").Value = 14687.75
$ws.Range("K73").Value = 14687.75
$ws.Range("M73").Value = -13751.75
$ws.Range("H113").Value = 302828.56
$ws.Range("I113").Value = 222960
$ws.Range("K113").Value = 222960
$ws.Range("M113").Value = -220790
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 474.75
$ws.Range("J16").Value = 350
$ws.Range("L16").Value = 350
$ws.Range("N16").Value = -690
$ws.Range("H22").Value = 2837.923
$ws.Range("I22").Value = 1299.8572
$ws.Range("J22").Value = 4632.3335
$ws.Range("K22").Value = 1299.8572
$ws.Range("L22").Value = 4632.3335
$ws.Range("M22").Value = -1004.8572
$ws.Range("N22").Value = -5222.3335
$ws.Range("H27").Value = 2837.923
$ws.Range("I27").Value = 1299.8572
$ws.Range("J27").Value = 4632.3335
$ws.Range("K27").Value = 1299.8572
$ws.Range("L27").Value = 4632.3335
$ws.Range("M27").Value = -1192.8572
$ws.Range("N27").Value = -4846.3335
$ws.Range("H46").Value = 8258.857
$ws.Range("I46").Value = 12761.75
$ws.Range("J46").Value = 5487.846
$ws.Range("K46").Value = 12761.75
$ws.Range("L46").Value = 5487.846
$ws.Range("M46").Value = -12573.75
$ws.Range("N46").Value = -5863.846
$ws.Range("H55").Value = 832.8125
$ws.Range("J55").Value = 1166.3334
$ws.Range("L55").Value = 1166.3334
$ws.Range("N55").Value = -1512.3334
$ws.Range("H82").Value = 2883.5334
$ws.Range("I82").Value = 2150
$ws.Range("J82").Value = 3066.9167
$ws.Range("K82").Value = 2150
$ws.Range("L82").Value = 3066.9167
$ws.Range("M82").Value = -1789
$ws.Range("N82").Value = -3788.9167
$ws.Range("H85").Value = 2883.5334
$ws.Range("I85").Value = 2150
$ws.Range("J85").Value = 3066.9167
$ws.Range("K85").Value = 2150
$ws.Range("L85").Value = 3066.9167
$ws.Range("M85").Value = -902
$ws.Range("N85").Value = -5562.9167
$ws.Range("H132").Value = 154374.12
$ws.Range("I132").Value = 175427.58
$ws.Range("K132").Value = 526282.74
$ws.Range("M132").Value = -523752.74
$ws.Range("H136").Value = 3826.4736
$ws.Range("I136").Value = 2646.6
$ws.Range("K136").Value = 7939.799999999999
$ws.Range("M136").Value = -5389.799999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1158.2174
$ws.Range("I113").Value = 1149.1578
$ws.Range("K113").Value = 3447.4734
$ws.Range("M113").Value = -1277.4734
$ws.Range("H132").Value = 33713.61
$ws.Range("I132").Value = 33713.61
$ws.Range("K132").Value = 101140.83
$ws.Range("M132").Value = -98610.83
